$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume-change (column E)
# values for rows 2-51, per the GitHub Actions data refresh.
# Numeric-looking price strings are entered with a leading apostrophe
# (Excel quote-prefix) so they are stored as text and keep their exact
# original formatting (trailing zeros, decimal places) instead of being
# auto-coerced into floating point numbers.

$ws.Range("D2").Value = '27.223.96'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.906.31'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''307.98'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '''0.5207'
$ws.Range("E7").Value = '  +1.00%  '
$ws.Range("D8").Value = '''0.3769'
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '''0.07276'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").Value = '''21.19'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '''0.9050'
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '''0.08447'
$ws.Range("E12").Value = '  +10.28%  '
$ws.Range("D13").Value = '1.913.81'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '''96.97'
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("D15").Value = '''5.300'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '''1.002'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '''0.000008674'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '''14.55'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '27.260.46'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").Value = '''5.097'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").Value = '2.158.34'
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = '''6.451'
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("D25").Value = '''2.329'
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").Value = '''146.85'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").Value = '''18.26'
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").Value = '''115.23'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").Value = '''4.829'
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").Value = '''4.914'
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '''0.09284'
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '''0.05080'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").Value = '''0.7964'
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("D35").Value = '''1.245'
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").Value = '''3.429'
$ws.Range("E36").Value = '  +4.49%  '
$ws.Range("D37").Value = '''2.954'
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").Value = '''0.5799'
$ws.Range("E38").Value = '  +3.27%  '
$ws.Range("D39").Value = '''2.589'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").Value = '''0.02009'
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '''9.073'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = '''6.620'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '''116.91'
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("D45").Value = '''0.1522'
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("D46").Value = '''0.4878'
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").Value = '''1.002'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '''10.16'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").Value = '''1.638'
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").Value = '''37.77'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '''64.17'
$ws.Range("E51").Value = '  +0.23%  '
